$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Periodo Mora" (column E) / "Valor Mora" (column F) table in B15:J20
# is being refreshed with the new account-statement data: periods are
# re-sorted into ascending order (2109, 2110, 2111, 2112, 2201) and the
# "Valor Mora" figure travels with its period.

$periodos = @("2109", "2110", "2111", "2112", "2201")
$valores  = @(48000, 48000, 48000, 48000, 40000)

for ($i = 0; $i -lt 5; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periodos[$i]
    $ws.Cells.Item($row, 6).Value = $valores[$i]
}
